$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value()
    $dVal = $dCell.Value()
    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
